$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("CapitalCall")

# Update Call Date (D2) to a new date and give it its own date-formatted style
$ws.Range("D2").Value = "10/22/2022"
$ws.Range("I2").Copy()
$ws.Range("D2").PasteSpecial(-4122)  # xlPasteFormats
$ws.Range("D2").NumberFormat = "mm-dd-yy"

# Due Date (E2) no longer derives from Call Date via formula; it now holds
# the previous Call Date value as a plain static date value
$ws.Range("E2").Value = "11/5/2022"
